$d = $word.ActiveDocument

# 1. Update the two "2025-01-25" dates (title date line and the "Date:" line in
#    the header table) to "2025-01-28". The unrelated "2025-01-24" date in the
#    version-history table must stay untouched, so use an exact, case-sensitive
#    whole match on "2025-01-25".
$d.Content.Find.Execute("2025-01-25", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-01-28", 2)

# 2. Add a new bullet item after the "interfaces and encodings for cloud
#    optimized data discovery and access" bullet, matching the same list
#    (numId 1002, ilvl 0).
$target = $d.Content
$target.Find.Execute("interfaces and encodings for cloud optimized data discovery and access",
                      $true, $false, $false, $false, $false,
                      $true, 1, $false, "", 0) | Out-Null

$target.Collapse(0)            # wdCollapseEnd
$target.InsertParagraphAfter()
$target.Collapse(0)            # wdCollapseEnd
$target.Move(4, 1) | Out-Null  # wdCharacter - step into the freshly created paragraph
$target.InsertAfter("data usage insights, weblogs for Global Caches and Global Discovery Catalogues")
